$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that no longer exist in the updated table (old rows 8 & 9) ---
$ws.Range("A8:E9").EntireRow.Delete()

# --- Re-write header row: columns reordered so Class moves from A to C,
#     and Area_Sq_Km / Proportion move from B/C to A/B ---
$ws.Cells.Item(1,1).Value = "Area_Sq_Km"
$ws.Cells.Item(1,2).Value = "Proportion"
$ws.Cells.Item(1,3).Value = "Class"
$ws.Cells.Item(1,4).Value = "Area_Exclusive"
$ws.Cells.Item(1,5).Value = "Area_Overlapped"

# --- Re-write the 6 data rows with the new column order & refreshed values ---
$ws.Cells.Item(2,1).Value = 30238.8931
$ws.Cells.Item(2,2).Value = 28.7932954632323
$ws.Cells.Item(2,3).Value = "Total"
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 0

$ws.Cells.Item(3,1).Value = 26284.8191
$ws.Cells.Item(3,2).Value = 25.0282495473987
$ws.Cells.Item(3,3).Value = "Nautra_2000"
$ws.Cells.Item(3,4).Value = 7274.0322
$ws.Cells.Item(3,5).Value = 19010.7869

$ws.Cells.Item(4,1).Value = 19305.923
$ws.Cells.Item(4,2).Value = 18.3829858881115
$ws.Cells.Item(4,3).Value = "Fredninger"
$ws.Cells.Item(4,4).Value = 203.7476
$ws.Cells.Item(4,5).Value = 19102.1754

$ws.Cells.Item(5,1).Value = 7157.3432
$ws.Cells.Item(5,2).Value = 6.81517993426012
$ws.Cells.Item(5,3).Value = "Havstrategi_standard"
$ws.Cells.Item(5,4).Value = 2264.7899
$ws.Cells.Item(5,5).Value = 4892.5533

$ws.Cells.Item(6,1).Value = 4300.9158
$ws.Cells.Item(6,2).Value = 4.09530662985426
$ws.Cells.Item(6,3).Value = "Havstrategi_streng"
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 4300.9158

$ws.Cells.Item(7,1).Value = 2378.0018
$ws.Cells.Item(7,2).Value = 2.26431927296632
$ws.Cells.Item(7,3).Value = "Natur_Vildt_Reservater"
$ws.Cells.Item(7,4).Value = 3.3514
$ws.Cells.Item(7,5).Value = 2374.6504
